$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1485.4615
$ws.Range("I62").Value = 1430.5
$ws.Range("J62").Value = 1668.6666
$ws.Range("K62").Value = 1430.5
$ws.Range("L62").Value = 1668.6666
$ws.Range("M62").Value = -806.5
$ws.Range("N62").Value = -2916.6666
$ws.Range("H65").Value = 1485.4615
$ws.Range("I65").Value = 1430.5
$ws.Range("J65").Value = 1668.6666
$ws.Range("K65").Value = 7152.5
$ws.Range("L65").Value = 8343.333000000001
$ws.Range("M65").Value = -4032.5
$ws.Range("N65").Value = -14583.333
$ws.Range("H76").Value = 4199.5
$ws.Range("I76").Value = 2626.5715
$ws.Range("J76").Value = 6401.6
$ws.Range("K76").Value = 2626.5715
$ws.Range("L76").Value = 6401.6
$ws.Range("M76").Value = -2311.5715
$ws.Range("N76").Value = -7031.6
$ws.Range("H79").Value = 4199.5
$ws.Range("I79").Value = 2626.5715
$ws.Range("J79").Value = 6401.6
$ws.Range("K79").Value = 2626.5715
$ws.Range("L79").Value = 6401.6
$ws.Range("M79").Value = -1534.5715
$ws.Range("N79").Value = -8585.6
$ws.Range("H80").Value = 630.0238000000001
$ws.Range("I80").Value = 452.5
$ws.Range("J80").Value = 739.2692
$ws.Range("K80").Value = 1357.5
$ws.Range("L80").Value = 2217.8076
$ws.Range("M80").Value = -359.5
$ws.Range("N80").Value = -4213.8076
$ws.Range("H83").Value = 630.0238000000001
$ws.Range("I83").Value = 452.5
$ws.Range("J83").Value = 739.2692
$ws.Range("K83").Value = 4072.5
$ws.Range("L83").Value = 6653.422799999999
$ws.Range("M83").Value = 919.5
$ws.Range("N83").Value = -16637.4228

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32260386
$ws.Range("I61").Value = 45456536
$ws.Range("J61").Value = 3125.5557
$ws.Range("K61").Value = 45456536
$ws.Range("L61").Value = 3125.5557
$ws.Range("M61").Value = -45456324
$ws.Range("N61").Value = -3549.5557
$ws.Range("H88").Value = 2830
$ws.Range("I88").Value = 2490
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 2490
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -2084
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 2830
$ws.Range("I91").Value = 2490
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 2490
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -1086
$ws.Range("N91").Value = -5808
$ws.Range("H122").Value = 1824
$ws.Range("I122").Value = 1527.75
$ws.Range("J122").Value = 2021.5
$ws.Range("K122").Value = 4583.25
$ws.Range("L122").Value = 6064.5
$ws.Range("M122").Value = -2133.25
$ws.Range("N122").Value = -10964.5
$ws.Range("H136").Value = 32260386
$ws.Range("I136").Value = 45456536
$ws.Range("J136").Value = 3125.5557
$ws.Range("K136").Value = 136369608
$ws.Range("L136").Value = 9376.667099999999
$ws.Range("M136").Value = -136367058
$ws.Range("N136").Value = -14476.6671

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 629159.3
$ws.Range("I86").Value = 3025.1667
$ws.Range("J86").Value = 1004839.8
$ws.Range("K86").Value = 3025.1667
$ws.Range("L86").Value = 1004839.8
$ws.Range("M86").Value = -1902.1667
$ws.Range("N86").Value = -1007085.8
$ws.Range("H89").Value = 629159.3
$ws.Range("I89").Value = 3025.1667
$ws.Range("J89").Value = 1004839.8
$ws.Range("K89").Value = 15125.8335
$ws.Range("L89").Value = 5024199
$ws.Range("M89").Value = -9509.833500000001
$ws.Range("N89").Value = -5035431
$ws.Range("H107").Value = 2049.6086
$ws.Range("I107").Value = 1858.7333
$ws.Range("J107").Value = 2407.5
$ws.Range("K107").Value = 1858.7333
$ws.Range("L107").Value = 2407.5
$ws.Range("M107").Value = 61.2666999999999
$ws.Range("N107").Value = -6247.5
$ws.Range("H134").Value = 3564.5264
$ws.Range("I134").Value = 3469.0908
$ws.Range("J134").Value = 3695.75
$ws.Range("K134").Value = 10407.2724
$ws.Range("L134").Value = 11087.25
$ws.Range("M134").Value = -7872.2724
$ws.Range("N134").Value = -16157.25

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 592.6799999999999
$ws.Range("I131").Value = 295.59616
$ws.Range("J131").Value = 914.5208
$ws.Range("K131").Value = 886.7884799999999
$ws.Range("L131").Value = 2743.5624
$ws.Range("M131").Value = 4153.21152
$ws.Range("N131").Value = -12823.5624

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 8000
$ws.Range("N57").Value = -9640
$ws.Range("H80").Value = 2974.5454
$ws.Range("I80").Value = 2143
$ws.Range("J80").Value = 3667.5
$ws.Range("K80").Value = 2143
$ws.Range("L80").Value = 3667.5
$ws.Range("M80").Value = -1145
$ws.Range("N80").Value = -5663.5
$ws.Range("H83").Value = 2974.5454
$ws.Range("I83").Value = 2143
$ws.Range("J83").Value = 3667.5
$ws.Range("K83").Value = 10715
$ws.Range("L83").Value = 18337.5
$ws.Range("M83").Value = -5723
$ws.Range("N83").Value = -28321.5
$ws.Range("H102").Value = 2281.5667
$ws.Range("I102").Value = 1982.6666
$ws.Range("J102").Value = 2979
$ws.Range("K102").Value = 1982.6666
$ws.Range("L102").Value = 2979
$ws.Range("M102").Value = -360.6666
$ws.Range("N102").Value = -6223
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1417.5
$ws.Range("I82").Value = 1228.6364
$ws.Range("J82").Value = 1714.2858
$ws.Range("K82").Value = 1228.6364
$ws.Range("L82").Value = 1714.2858
$ws.Range("M82").Value = -867.6364000000001
$ws.Range("N82").Value = -2436.2858
$ws.Range("H85").Value = 1417.5
$ws.Range("I85").Value = 1228.6364
$ws.Range("J85").Value = 1714.2858
$ws.Range("K85").Value = 1228.6364
$ws.Range("L85").Value = 1714.2858
$ws.Range("M85").Value = 19.36359999999991
$ws.Range("N85").Value = -4210.2858
$ws.Range("H100").Value = 1301.2
$ws.Range("J100").Value = 1600
$ws.Range("L100").Value = 1600
$ws.Range("N100").Value = -2682
$ws.Range("H132").Value = 22366.143
$ws.Range("I132").Value = 34151
$ws.Range("J132").Value = 6653
$ws.Range("K132").Value = 102453
$ws.Range("L132").Value = 19959
$ws.Range("M132").Value = -99923
$ws.Range("N132").Value = -25019
$ws.Range("H139").Value = 54000
$ws.Range("J139").Value = 54000
$ws.Range("L139").Value = 54000
$ws.Range("N139").Value = -64280

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45700
$ws.Range("J46").Value = 45700
$ws.Range("L46").Value = 45700
$ws.Range("N46").Value = -46162
$ws.Range("H81").Value = 846.1
$ws.Range("I81").Value = 765.8570999999999
$ws.Range("J81").Value = 1033.3334
$ws.Range("K81").Value = 1531.7142
$ws.Range("L81").Value = 2066.6668
$ws.Range("M81").Value = -470.7141999999999
$ws.Range("N81").Value = -4188.6668
$ws.Range("H84").Value = 846.1
$ws.Range("I84").Value = 765.8570999999999
$ws.Range("J84").Value = 1033.3334
$ws.Range("K84").Value = 7658.571
$ws.Range("L84").Value = 10333.334
$ws.Range("M84").Value = -2354.571
$ws.Range("N84").Value = -20941.334
$ws.Range("H86").Value = 17183.223
$ws.Range("J86").Value = 17183.223
$ws.Range("L86").Value = 17183.223
$ws.Range("N86").Value = -19429.223
$ws.Range("H89").Value = 17183.223
$ws.Range("J89").Value = 17183.223
$ws.Range("L89").Value = 85916.11500000001
$ws.Range("N89").Value = -97148.11500000001
$ws.Range("H96").Value = 1860
$ws.Range("J96").Value = 2433.3333
$ws.Range("L96").Value = 2433.3333
$ws.Range("N96").Value = -5179.3333
$ws.Range("H132").Value = 7350.8696
$ws.Range("I132").Value = 9879.385
$ws.Range("J132").Value = 4063.8
$ws.Range("K132").Value = 29638.155
$ws.Range("L132").Value = 12191.4
$ws.Range("M132").Value = -27108.155
$ws.Range("N132").Value = -17251.4
$ws.Range("H134").Value = 45700
$ws.Range("J134").Value = 45700
$ws.Range("L134").Value = 137100
$ws.Range("N134").Value = -142170
$ws.Range("H136").Value = 3118.652
$ws.Range("I136").Value = 4245.107
$ws.Range("J136").Value = 1366.3889
$ws.Range("K136").Value = 12735.321
$ws.Range("L136").Value = 4099.1667
$ws.Range("M136").Value = -10185.321
$ws.Range("N136").Value = -9199.1667
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
